$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 3 blank rows right after the first data row (old row 16, KAREN)
#    so the data block grows from 9 rows (16-24) to 12 rows (16-27).
# ---------------------------------------------------------------------------
$ws.Range("A17:A19").EntireRow.Insert()

# Fix up the formatting of the 3 freshly-inserted blank rows so they match
# the "normal" data-row style (copy formats from row 16).
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J19").PasteSpecial(-4122)
$ws.Range("B16:J16").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Overwrite the whole data block (rows 16-27) with the final values.
#    Rows 16-26 all share the "normal" row style; row 27 keeps the
#    special "last row" style that already lives on the last row of the
#    block (it was originally row 24, now pushed down to row 27 by the
#    insert above, so its distinctive bottom-border formatting is already
#    in place - we only need to update its values).
# ---------------------------------------------------------------------------

# Row 16: WILMER JOSE PUELLO GONZALEZ
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1143329024"
$ws.Range("D16").Value = "WILMER JOSE PUELLO GONZALEZ"
$ws.Range("E16").Value = "2506"
$ws.Range("F16").Value = 17067
$ws.Range("G16").Value = 1600000

# Row 17: IRAYDA LAVOE ALCALA RIOS (new)
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1047486394"
$ws.Range("D17").Value = "IRAYDA LAVOE ALCALA RIOS"
$ws.Range("E17").Value = "2507"
$ws.Range("F17").Value = 5694
$ws.Range("G17").Value = 1423500

# Row 18: OSNAIDER ENRIQUE RIVERO RIVERO (new)
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1050951430"
$ws.Range("D18").Value = "OSNAIDER ENRIQUE RIVERO RIVERO"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 1898
$ws.Range("G18").Value = 1423500

# Row 19: ERIS ANDRES ARELLANO CABARCAS (new)
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1007264496"
$ws.Range("D19").Value = "ERIS ANDRES ARELLANO CABARCAS"
$ws.Range("E19").Value = "2507"
$ws.Range("F19").Value = 1898
$ws.Range("G19").Value = 1423500

# Row 20: KAREN SILGADO AYALA
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1128052700"
$ws.Range("D20").Value = "KAREN SILGADO AYALA"
$ws.Range("E20").Value = "1803"
$ws.Range("F20").Value = 31249
$ws.Range("G20").Value = 781242

# Row 21: EMIRO RAFAEL MARTINEZ BENITEZ - 2008
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1068659636"
$ws.Range("D21").Value = "EMIRO RAFAEL MARTINEZ BENITEZ"
$ws.Range("E21").Value = "2008"
$ws.Range("F21").Value = 35112
$ws.Range("G21").Value = 1000000

# Row 22: EMIRO RAFAEL MARTINEZ BENITEZ - 2007
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1068659636"
$ws.Range("D22").Value = "EMIRO RAFAEL MARTINEZ BENITEZ"
$ws.Range("E22").Value = "2007"
$ws.Range("F22").Value = 35112
$ws.Range("G22").Value = 1000000

# Row 23: EMIRO RAFAEL MARTINEZ BENITEZ - 2006
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1068659636"
$ws.Range("D23").Value = "EMIRO RAFAEL MARTINEZ BENITEZ"
$ws.Range("E23").Value = "2006"
$ws.Range("F23").Value = 35112
$ws.Range("G23").Value = 1000000

# Row 24: EMIRO RAFAEL MARTINEZ BENITEZ - 2005
$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1068659636"
$ws.Range("D24").Value = "EMIRO RAFAEL MARTINEZ BENITEZ"
$ws.Range("E24").Value = "2005"
$ws.Range("F24").Value = 35112
$ws.Range("G24").Value = 1000000

# Row 25: EMIRO RAFAEL MARTINEZ BENITEZ - 2004
$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1068659636"
$ws.Range("D25").Value = "EMIRO RAFAEL MARTINEZ BENITEZ"
$ws.Range("E25").Value = "2004"
$ws.Range("F25").Value = 35112
$ws.Range("G25").Value = 1000000

# Row 26: EMIRO RAFAEL MARTINEZ BENITEZ - 2003
$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1068659636"
$ws.Range("D26").Value = "EMIRO RAFAEL MARTINEZ BENITEZ"
$ws.Range("E26").Value = "2003"
$ws.Range("F26").Value = 10534
$ws.Range("G26").Value = 1000000

# Row 27: LUIS ANGEL ACOSTA HERNANDEZ (last row of table - keeps the
# distinctive bottom-border styling that already lives on this row).
$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "1044934463"
$ws.Range("D27").Value = "LUIS ANGEL ACOSTA HERNANDEZ"
$ws.Range("E27").Value = "2301"
$ws.Range("F27").Value = 21654
$ws.Range("G27").Value = 1160000

# ---------------------------------------------------------------------------
# 3. Update the worker / period counts and the total overdue amount.
# ---------------------------------------------------------------------------
$ws.Range("C13").Value = 7
$ws.Range("F13").Value = 10
$ws.Range("E11").Value = 265554
